$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B4").Value = '${from.toString("YYYY.MM.dd HH:mm:ss")+" - "+to.toString("YYYY.MM.dd HH:mm:ss")}'
